# Weekly CompStat report refresh — new crime data collected.
# Updates the report header (volume/week-of dates) and the precinct
# crime-complaint statistics table (rows 15-28) to the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helpers -----------------------------------------------------------
# Plain numeric overwrite; cell keeps its existing style/number format.
function Set-Num {
    param($ws, $addr, $val)
    $ws.Range($addr).Value = $val
}

# Some "% Chg" cells flip between a real number and the literal marker
# text "0" / "***.*" (used when the prior-period base is zero). Those
# marker cells use a distinct right-aligned text style (shared with the
# header row). Converting a text-marker cell back to a number: copy the
# format from a stable numeric donor cell, then write the number.
function Set-NumWithDonor {
    param($ws, $addr, $val, $donor)
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $val
}

# Converting a numeric cell to the literal marker text: copy format AND
# value from a donor cell that already holds the right marker text/style
# (-4122 = xlPasteFormats, -4163 = xlPasteValues), so the result keeps
# the exact shared-string-backed text cell type used elsewhere.
function Set-TextFromDonor {
    param($ws, $addr, $donor)
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# Stable donor cells (outside the edited row range, untouched by this
# week's refresh) used purely as format/value sources above:
#   C14 -> style 13 text "0"      E14 -> style 13 text "***.*"
#   F31 -> style 15 plain number  L14 -> style 14 plain number

# --- header: volume number + week-of date range -------------------------
$ws.Range("A8").Value = "Volume 32   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/17/2025  Through  3/23/2025"

# --- crime complaint table (rows 15-28) ---------------------------------
Set-TextFromDonor $ws "G15" "C14"
Set-TextFromDonor $ws "H15" "E14"
Set-Num $ws "L15" -25
Set-TextFromDonor $ws "C16" "C14"
Set-Num $ws "D16" 2
Set-Num $ws "E16" -100
Set-Num $ws "G16" 4
Set-Num $ws "H16" 50
Set-Num $ws "J16" 18
Set-Num $ws "K16" -5.555555555555
Set-Num $ws "L16" -26.086956521739
Set-Num $ws "N16" -87.218045112782
Set-Num $ws "D17" 2
Set-Num $ws "E17" 0
Set-Num $ws "F17" 6
Set-Num $ws "G17" 9
Set-Num $ws "H17" -33.333333333333
Set-Num $ws "I17" 16
Set-Num $ws "J17" 25
Set-Num $ws "K17" -36
Set-Num $ws "L17" -27.272727272727
Set-Num $ws "M17" 0
Set-Num $ws "N17" -15.789473684210
Set-Num $ws "C18" 1
Set-Num $ws "D18" 8
Set-Num $ws "E18" -87.5
Set-Num $ws "F18" 8
Set-Num $ws "G18" 23
Set-Num $ws "H18" -65.217391304347
Set-Num $ws "I18" 38
Set-Num $ws "J18" 47
Set-Num $ws "K18" -19.148936170212
Set-Num $ws "L18" -5
Set-Num $ws "M18" 11.764705882352
Set-Num $ws "N18" -80
Set-Num $ws "C19" 16
Set-Num $ws "D19" 12
Set-Num $ws "E19" 33.333333333333
Set-Num $ws "F19" 51
Set-Num $ws "H19" -1.923076923076
Set-Num $ws "I19" 139
Set-Num $ws "J19" 163
Set-Num $ws "K19" -14.723926380368
Set-Num $ws "L19" -7.333333333333
Set-Num $ws "M19" -15.243902439024
Set-Num $ws "N19" -69.782608695652
Set-NumWithDonor $ws "C20" 1 "F31"
Set-TextFromDonor $ws "D20" "C14"
Set-TextFromDonor $ws "E20" "E14"
Set-Num $ws "F20" 2
Set-Num $ws "G20" 3
Set-Num $ws "H20" -33.333333333333
Set-Num $ws "I20" 6
Set-Num $ws "K20" -14.285714285714
Set-Num $ws "L20" -64.705882352941
Set-Num $ws "M20" 50
Set-Num $ws "N20" -97.794117647058
Set-Num $ws "C21" 20
Set-Num $ws "D21" 24
Set-Num $ws "E21" -16.666666666666
Set-Num $ws "F21" 75
Set-Num $ws "G21" 91
Set-Num $ws "H21" -17.582417582417
Set-Num $ws "I21" 219
Set-Num $ws "J21" 262
Set-Num $ws "K21" -16.412213740458
Set-Num $ws "L21" -14.785992217898
Set-Num $ws "M21" -8.368200836820
Set-Num $ws "N21" -79.722222222222
Set-Num $ws "F22" 1
Set-Num $ws "H22" -66.666666666666
Set-Num $ws "J22" 9
Set-Num $ws "K22" -44.444444444444
Set-Num $ws "D23" 1
Set-Num $ws "G23" 4
Set-Num $ws "H23" -25
Set-Num $ws "J23" 10
Set-Num $ws "K23" -30
Set-Num $ws "L23" -46.153846153846
Set-Num $ws "C24" 17
Set-Num $ws "D24" 12
Set-Num $ws "E24" 41.666666666666
Set-Num $ws "F24" 99
Set-Num $ws "G24" 76
Set-Num $ws "H24" 30.263157894736
Set-Num $ws "I24" 273
Set-Num $ws "J24" 253
Set-Num $ws "K24" 7.905138339920
Set-Num $ws "L24" 20.264317180616
Set-Num $ws "M24" 34.482758620689
Set-Num $ws "C25" 15
Set-Num $ws "D25" 10
Set-Num $ws "E25" 50
Set-Num $ws "F25" 78
Set-Num $ws "G25" 59
Set-Num $ws "H25" 32.203389830508
Set-Num $ws "I25" 230
Set-Num $ws "J25" 195
Set-Num $ws "K25" 17.948717948717
Set-Num $ws "L25" 20.418848167539
Set-Num $ws "D26" 8
Set-Num $ws "E26" -62.5
Set-Num $ws "F26" 16
Set-Num $ws "G26" 18
Set-Num $ws "H26" -11.111111111111
Set-Num $ws "I26" 45
Set-Num $ws "J26" 51
Set-Num $ws "K26" -11.764705882352
Set-Num $ws "L26" -10
Set-Num $ws "M26" -15.094339622641
Set-Num $ws "G27" 1
Set-Num $ws "H27" 100
Set-Num $ws "L27" 0
Set-TextFromDonor $ws "C28" "C14"
Set-NumWithDonor $ws "D28" 1 "F31"
Set-NumWithDonor $ws "E28" -100 "L14"
Set-NumWithDonor $ws "G28" 1 "F31"
Set-NumWithDonor $ws "H28" 300 "L14"
Set-Num $ws "J28" 3
Set-Num $ws "K28" 166.666666666667
